$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clusters")

$ws.Range("C4").Value  = "property, corporation, interest, stock, foreign, business, item, basis, without, loss"
$ws.Range("C5").Value  = "purposes, term, case, general, qualified, described, apply, means, provided, defined"
$ws.Range("C6").Value  = "shall, subsection, paragraph, person, subparagraph, stat, states, united, percent, trust"
$ws.Range("C7").Value  = "year, taxable, amount, tax, respect, income, plan, treated, period, date"
$ws.Range("C9").Value  = "secretary, may, regulations, return, required, information"
$ws.Range("C10").Value = "shall, subsection, amount, paragraph, respect, income, property, corporation, plan, person"
$ws.Range("C11").Value = "purposes, term, case, general, described, qualified, treated, apply, made, means"
$ws.Range("C12").Value = "tax, stat, states, foreign, credit, certain, item, part, deduction, act"
$ws.Range("C13").Value = "year, taxable, date, period, percent, business, years, loss, calendar, beginning"
